# Apply updated cryptos list (prices in column D, 1h volume % in column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to stay a text string (not auto-parsed as a number/date)
    # while keeping its original style (no $ s= attribute change).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "24.666.02"
$ws.Range("E2").Value = "  +0.51%  "
Set-TextValue $ws.Range("D3") "1.692.47"
$ws.Range("E3").Value = "  +0.06%  "
Set-TextValue $ws.Range("D4") "1.004"
$ws.Range("E4").Value = "  +0.24%  "
Set-TextValue $ws.Range("D5") "316.76"
$ws.Range("E5").Value = "  +1.39%  "
Set-TextValue $ws.Range("D6") "1.002"
$ws.Range("E6").Value = "  +0.22%  "
Set-TextValue $ws.Range("D7") "0.3953"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("E8").Value = "  +1.04%  "
Set-TextValue $ws.Range("D9") "1.492"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("E10").Value = "  +0.29%  "
Set-TextValue $ws.Range("D11") "52.07"
$ws.Range("E11").Value = "  -2.56%  "
Set-TextValue $ws.Range("D12") "0.08886"
$ws.Range("E12").Value = "  +1.51%  "
Set-TextValue $ws.Range("D13") "7.250"
$ws.Range("E13").Value = "  -0.63%  "
Set-TextValue $ws.Range("D14") "23.64"
$ws.Range("E14").Value = "  +2.10%  "
Set-TextValue $ws.Range("D15") "8.070"
$ws.Range("E15").Value = "  +7.03%  "
Set-TextValue $ws.Range("D16") "0.00001321"
$ws.Range("E16").Value = "  +0.38%  "
Set-TextValue $ws.Range("D17") "1.696.81"
$ws.Range("E17").Value = "  +0.14%  "
Set-TextValue $ws.Range("D18") "99.94"
$ws.Range("E18").Value = "  -0.26%  "
Set-TextValue $ws.Range("D19") "0.07029"
$ws.Range("E19").Value = "  -0.60%  "
Set-TextValue $ws.Range("D20") "19.62"
$ws.Range("E20").Value = "  +0.94%  "
Set-TextValue $ws.Range("D21") "7.005"
$ws.Range("E21").Value = "  +4.68%  "
Set-TextValue $ws.Range("D22") "1.007"
$ws.Range("E22").Value = "  +0.53%  "
Set-TextValue $ws.Range("D23") "14.34"
$ws.Range("E23").Value = "  +1.59%  "
Set-TextValue $ws.Range("D24") "24.647.63"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("E25").Value = "  +6.39%  "
$ws.Range("E26").Value = "  +2.21%  "
Set-TextValue $ws.Range("D27") "22.74"
$ws.Range("E27").Value = "  +2.09%  "
Set-TextValue $ws.Range("D28") "162.29"
Set-TextValue $ws.Range("D29") "135.97"
$ws.Range("E29").Value = "  +1.93%  "
Set-TextValue $ws.Range("D30") "5.190"
$ws.Range("E30").Value = "  +0.76%  "
Set-TextValue $ws.Range("D31") "7.566"
$ws.Range("E31").Value = "  +0.65%  "
Set-TextValue $ws.Range("D32") "0.08611"
$ws.Range("E32").Value = "  -0.23%  "
Set-TextValue $ws.Range("D33") "1.055"
$ws.Range("E33").Value = "  -3.14%  "
Set-TextValue $ws.Range("D34") "7.055"
$ws.Range("E34").Value = "  -3.62%  "
Set-TextValue $ws.Range("D35") "11.37"
$ws.Range("E35").Value = "  +3.42%  "
Set-TextValue $ws.Range("D36") "0.2734"
$ws.Range("E36").Value = "  +0.91%  "
Set-TextValue $ws.Range("D37") "1.886"
$ws.Range("E37").Value = "  -4.09%  "
Set-TextValue $ws.Range("D38") "14.46"
$ws.Range("E38").Value = "  -1.59%  "
Set-TextValue $ws.Range("D39") "0.09228"
$ws.Range("E39").Value = "  +2.85%  "
Set-TextValue $ws.Range("D40") "0.02725"
$ws.Range("E40").Value = "  -0.95%  "
Set-TextValue $ws.Range("D41") "1.471"
$ws.Range("E41").Value = "  -0.07%  "
Set-TextValue $ws.Range("D42") "0.7662"
$ws.Range("E42").Value = "  +0.47%  "
Set-TextValue $ws.Range("D43") "16.20"
$ws.Range("E43").Value = "  +5.21%  "
Set-TextValue $ws.Range("D44") "2.599"
$ws.Range("E44").Value = "  +6.36%  "
Set-TextValue $ws.Range("D45") "0.7159"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("E47").Value = "  +0.29%  "
Set-TextValue $ws.Range("D48") "140.27"
$ws.Range("E48").Value = "  +0.08%  "
Set-TextValue $ws.Range("D49") "1.320"
$ws.Range("E49").Value = "  +2.20%  "
Set-TextValue $ws.Range("D50") "90.91"
$ws.Range("E50").Value = "  +5.48%  "
Set-TextValue $ws.Range("D51") "0.07987"
$ws.Range("E51").Value = "  -0.06%  "
